$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.750.59'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '2.493.11'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''533.56'
$ws.Range("E5").Value = '  +5.48%  '
$ws.Range("D6").Value = '''134.12'
$ws.Range("E6").Value = '  +3.70%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +2.90%  '
$ws.Range("D9").Value = '2.517.05'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").Value = '''0.0995'
$ws.Range("E10").Value = '  +4.24%  '
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").Value = '2.941.43'
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("D15").Value = '58.640.05'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '''22.34'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").Value = '2.504.82'
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").Value = '''10.65'
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").Value = '''4.25'
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("D21").Value = '''320.90'
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("D22").Value = '''6.21'
$ws.Range("E22").Value = '  +9.49%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '''65.51'
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("D25").Value = '''0.411'
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("D26").Value = '''0.996'
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '''7.50'
$ws.Range("E28").Value = '  +3.70%  '
$ws.Range("D29").Value = '''0.0₃0762'
$ws.Range("E29").Value = '  +5.21%  '
$ws.Range("D30").Value = '''172.65'
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("E31").Value = '  +5.28%  '
$ws.Range("E32").Value = '  +4.52%  '
$ws.Range("D33").Value = '''6.31'
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '''0.994'
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").Value = '''18.16'
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("D38").Value = '''3.95'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("E39").Value = '  +4.16%  '
$ws.Range("D40").Value = '''36.76'
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("D41").Value = '''0.822'
$ws.Range("E41").Value = '  +8.23%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''5.18'
$ws.Range("E42").Value = '  +4.57%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''3.48'
$ws.Range("E43").Value = '  +3.28%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''277.07'
$ws.Range("E44").Value = '  +2.30%  '
$ws.Range("D45").Value = '''131.28'
$ws.Range("E45").Value = '  +9.73%  '
$ws.Range("D46").Value = '''0.592'
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D47").Value = '''0.0935'
$ws.Range("E47").Value = '  +2.67%  '
$ws.Range("D48").Value = '''0.0510'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("E49").Value = '  +4.98%  '
$ws.Range("D50").Value = '''17.06'
$ws.Range("E50").Value = '  +2.68%  '
$ws.Range("D51").Value = '1.753.09'
$ws.Range("E51").Value = '  +3.23%  '
